$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "40*40"
$ws.Range("K21").Value = "VGGNet + LSTM"
$ws.Range("C21").Value = "VGG16 + LSTM"
$ws.Range("D21").Value = 0.1
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = 0.5
$ws.Range("G21").Value = 0.2
$ws.Range("H21").Value = 0.1
$ws.Range("J21").Value = "0-.25"

$ws.Range("A16:K16").Interior.Color = 5296274

$ws.Range("G25").Select()
